# MSME Country Indicators - Egypt, Arab Rep. Summary
# Rename the "Data" sheet to "Summary", insert a new "Source Type" header
# row, add a missing MSME employment figure, and append the AFDB source
# citation block at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet -----------------------------------------------------
$ws.Name = "Summary"

# --- Make room for the new "Source Type" header ------------------------
# Old rows 5-11 (the data table + source line) need to end up at rows
# 11-17, so insert 6 blank rows starting at row 5.
$ws.Range("5:10").Insert()

# --- New row 9: bold + underlined "Source Type" sub-header -------------
$ws.Range("A9").Value = "Source Type: SME Associations (Most Widely Used)"
$ws.Range("A9").Font.Bold = $true
$ws.Range("A9").Font.Underline = $true

# --- New MSME employment (% of total) figure in D14 --------------------
# Force text formatting so the value is stored the same way as its
# neighbouring cells (e.g. "20.5") rather than as a numeric value.
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.8"

# --- New source / citation block at the bottom of the sheet ------------
$ws.Range("A23").Value = "AFDB"
$ws.Range("A24").Value = 'African Development Bank, "Egypt Private Sector Country Profile", 2009, p. 115, 113. Available at http://www.afdb.org/fileadmin/uploads/afdb/Documents/Project-and-Operations/Brochure%20Egypt%20Anglais.pdf'

# --- Reassert formatting on cells whose rows were shifted ---------------
# (the row insert above resets the visual formatting of any cell it
# moves, so the "name"/"title"/"source" looks need to be re-applied.)
$ws.Range("A1").Font.Size = 18

$ws.Range("A3").Font.Bold = $true
$ws.Range("B11:D11").Font.Bold = $true
$ws.Range("A12").Font.Bold = $true
$ws.Range("A13").Font.Bold = $true
$ws.Range("A14").Font.Bold = $true
$ws.Range("A15").Font.Bold = $true
$ws.Range("A16").Font.Bold = $true
$ws.Range("A23").Font.Bold = $true

$ws.Range("A17").Font.Italic = $true
$ws.Range("A24").Font.Italic = $true
